# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (E) rows for the two trabajadores (DAVID JOSE RUIZ
# PUERTA / 1047456907 and JAIDER DARIO ZAPATEIRO SALGADO / 1143387469)
# are re-sequenced so the periods interleave worker-by-worker in
# chronological order (1812, 1901..1907) instead of being grouped by
# worker with a descending period order. The "Valor Mora" (F) amounts
# move together with their period (15625 stays tied to period 1907,
# 31249 to every other period).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047456907"
$ws.Range("D16").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E16").Value = "1812"
$ws.Range("F16").Value = 31249

$ws.Range("C17").Value = "1143387469"
$ws.Range("D17").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E17").Value = "1812"
$ws.Range("F17").Value = 31249

$ws.Range("C18").Value = "1047456907"
$ws.Range("D18").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E18").Value = "1901"
$ws.Range("F18").Value = 31249

$ws.Range("C19").Value = "1143387469"
$ws.Range("D19").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E19").Value = "1901"
$ws.Range("F19").Value = 31249

$ws.Range("C20").Value = "1047456907"
$ws.Range("D20").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E20").Value = "1902"
$ws.Range("F20").Value = 31249

$ws.Range("C21").Value = "1143387469"
$ws.Range("D21").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E21").Value = "1902"
$ws.Range("F21").Value = 31249

$ws.Range("C22").Value = "1047456907"
$ws.Range("D22").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E22").Value = "1903"
$ws.Range("F22").Value = 31249

$ws.Range("C23").Value = "1143387469"
$ws.Range("D23").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E23").Value = "1903"
$ws.Range("F23").Value = 31249

$ws.Range("C24").Value = "1047456907"
$ws.Range("D24").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E24").Value = "1904"
$ws.Range("F24").Value = 31249

$ws.Range("C25").Value = "1143387469"
$ws.Range("D25").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E25").Value = "1904"
$ws.Range("F25").Value = 31249

$ws.Range("C26").Value = "1047456907"
$ws.Range("D26").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E26").Value = "1905"
$ws.Range("F26").Value = 31249

$ws.Range("C27").Value = "1143387469"
$ws.Range("D27").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E27").Value = "1905"
$ws.Range("F27").Value = 31249

$ws.Range("C28").Value = "1047456907"
$ws.Range("D28").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E28").Value = "1906"
$ws.Range("F28").Value = 31249

$ws.Range("C29").Value = "1143387469"
$ws.Range("D29").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E29").Value = "1906"
$ws.Range("F29").Value = 31249

$ws.Range("C30").Value = "1047456907"
$ws.Range("D30").Value = "DAVID JOSE RUIZ PUERTA"
$ws.Range("E30").Value = "1907"
$ws.Range("F30").Value = 15625

$ws.Range("C31").Value = "1143387469"
$ws.Range("D31").Value = "JAIDER DARIO ZAPATEIRO SALGADO"
$ws.Range("E31").Value = "1907"
$ws.Range("F31").Value = 15625
